$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 79.42851574123137
    3  = 123.2952574732174
    4  = 35.57619318654503
    5  = 68.98063261057932
    6  = 71.93161907612469
    7  = 41.25751180954379
    8  = 85.5818596394662
    9  = 66.29500823011963
    10 = 90.65587845786484
    11 = 59.38041122973839
    12 = 67.24773137727452
    13 = 49.04403004832292
    14 = 46.67766993614866
    15 = 100.9536687017665
    16 = 52.01816211634516
    17 = 76.86563734451821
    18 = 54.18419287058687
    19 = 59.51531127249995
    20 = 73.15315572356192
    21 = 39.02044398478295
    22 = 70.70041955415016
    23 = 146.4679487280519
    24 = 88.98555492173882
    25 = 50.98030365658955
    26 = 73.67547077190636
    27 = 40.73555247053601
    28 = 46.71267873852356
    29 = 76.35929319148482
    30 = 35.30276851946299
    31 = 59.84841139159953
    32 = 77.4295813193125
    33 = 148.0919216512835
    34 = 152.6006981941599
    35 = 45.8292545256063
    36 = 214.3355849597693
    37 = 42.0482411963317
    38 = 67.36765891228413
    39 = 41.09569937623688
    40 = 45.3343464284967
    41 = 242.7265399639927
    42 = 35.40435868230194
    43 = 65.13417548341353
    44 = 106.6127349756764
    45 = 60.1188685680356
    46 = 135.8350386123793
    47 = 53.86402500069548
    48 = 45.89301057888606
    49 = 68.59910919115021
    50 = 141.9424880444684
    51 = 50.41423878759534
    52 = 49.09888572951384
}

foreach ($row in $values.Keys) {
    $ws.Range("T$row").Value = $values[$row]
}
